$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Objetivos (row 10/11 area) ---
$ws.Range("B10").Value = "Possibilitar aos alunos a realização de trabalho de síntese e integração dos conhecimentos adquiridos ao longo do curso,conforme projeto aprovado na disciplina de Trabalho de Graduação em Engenharia de Produção I."
$ws.Range("C10").Value = "Possibilitar aos alunos a realização de trabalho de síntese e integração dos conhecimentos adquiridos ao longo do curso,conforme projeto aprovado na disciplina de Trabalho de Graduação em Engenharia de Produção I."

# --- Docentes responsáveis now gets two extra rows inserted (rows 13/14) ---
$ws.Rows("13:14").Insert()

$ws.Range("A13").Value = ""
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows("13").RowHeight = 15

$ws.Range("B14").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C14").Value = "1285870 - Marcos Villela Barcza"
$ws.Rows("14").RowHeight = 15

# --- Programa resumido content (row 15, was row 13) ---
$ws.Range("B15").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se numprojeto de engenharia de produção."
$ws.Range("C15").Value = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se numprojeto de engenharia de produção."
$ws.Rows("15").RowHeight = 60

$ws.Rows("16").RowHeight = 60

# --- Programa (row 17, was row 15's A label, now gets content) ---
$ws.Range("B17").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) arevisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7)as conclusões e recomendações para trabalhos futuros e (8) referências bibliográficas."
$ws.Range("C17").Value = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) arevisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7)as conclusões e recomendações para trabalhos futuros e (8) referências bibliográficas."

# --- clear the stale old content that used to sit in rows that are now just labels ---
$ws.Range("B18:C18").ClearContents()
$ws.Range("B19:C19").ClearContents()

# --- Método content shifts down one row due to the earlier insert ---
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final."
$ws.Range("C20").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final."
$ws.Rows("21").RowHeight = 60

$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Uma única prova perante uma banca com 3 examinadores. A nota da disciplina será decidida pelos docentes da banca."
$ws.Range("C21").Value = "Uma única prova perante uma banca com 3 examinadores. A nota da disciplina será decidida pelos docentes da banca."

$ws.Rows("22").RowHeight = 60
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Reapresentação do trabalho modificado para nova avaliação."
$ws.Range("C22").Value = "Reapresentação do trabalho modificado para nova avaliação."

$ws.Rows("23").Insert()
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "Recomendada pelo orientador"
$ws.Range("C23").Value = "Recomendada pelo orientador"
$ws.Rows("23").RowHeight = 120

$ws.Range("A24").Value = "Requisitos:"

$ws.Range("B25").Value = "LOQ4228 -  Trabalho de Graduação em Engenharia de Produção I  (Requisito)`n"
$ws.Range("C25").Value = "LOQ4228 -  Trabalho de Graduação em Engenharia de Produção I  (Requisito)`n"
$ws.Rows("25").RowHeight = 30
